# Manual annotation: add "User Satisfication" (D) and "Relevance" (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
$ws.Range("D1").Value = "User Satisfication"
$ws.Range("E1").Value = "Relevance"

# E1 ("Relevance") -> bold, theme-colored font, no border (matches font used for bold
# default-styled cells in this workbook)
$ws.Range("E1").Font.Bold = $true

# D1 ("User Satisfication") -> bold Arial font (same family as the A1:C1 header font)
# with a thin left+right border and centered/top alignment, built by copying the
# existing header format (A1) and then stripping the top/bottom border edges.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Borders.Item(8).LineStyle = -4142
$ws.Range("D1").Borders.Item(9).LineStyle = -4142
$ws.Application.CutCopyMode = $false

# --- Column widths for the two new columns ---
$ws.Columns.Item(4).ColumnWidth = 23
$ws.Columns.Item(5).ColumnWidth = 27.25

# --- Data values for rows 2-26 ---
$values = @{
    2  = 1; 3  = 1; 4  = 1; 5  = 1; 6  = 1; 7  = 1; 8  = 1; 9  = 1
    10 = 0; 11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1
    18 = 1; 19 = 1; 20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1
}
foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 4).Value = $values[$r]
    $ws.Cells.Item($r, 5).Value = $values[$r]
}

# Row 26: User Satisfication = 0, Relevance = 1
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = 1

# --- View / selection / scroll position ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("E26").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
